$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Mark the first data row (CBSA 40900) as an "Example" row: fill A2:C2 red
#    and change its Comment from "Preferred" to "Example".
$ws.Range("A2:C2").Interior.Color = 255
$ws.Range("D2").Value = "Example"

# 2. The placeholder/error row (CBSA "XXXXX" / #N/A), originally row 7, is
#    moved down to the very bottom of the table and its Comment is changed
#    from "Preferred" to "NA". Deleting row 7 shifts every row below it up
#    by one, which is exactly what is needed before re-adding the row's
#    data at the new last row (152).
$ws.Rows.Item(7).Delete()

# Reset the AutoFilter while the table still ends at row 151 (i.e. before
# the moved row is written back), so the filter range covers the normal
# data block only (A1:D151), matching the post-edit filter/table extent.
$ws.AutoFilterMode = $false
$ws.Range("A1:D151").AutoFilter()

# Write the moved row back at the bottom of the sheet.
$ws.Range("A152").Value = "XXXXX"
$ws.Range("B152").Value = 45439
$ws.Range("C152").Formula = "=NA()"
$ws.Range("D152").Value = "NA"

# Keep the hidden _FilterDatabase defined name in sync with the new
# AutoFilter range.
$wb.Names.Item("trip_hh_cbsa_list!_FilterDatabase").RefersTo = "=trip_hh_cbsa_list!`$A`$1:`$D`$151"

# 3. Leave the selection where the editor left off.
$ws.Range("B5").Select()
